$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Destination")

$ws.Range("B1").Value = "D:/Red-Nord/DISC T/REGISTRU AUTORIZATII"
$ws.Range("B2").Value = "D:/Red-Nord/DISC T/DECONECTARI ZILNICE"
$ws.Range("B3").Value = "D:/Red-Nord/DISC T/RAPORT PDJT"
$ws.Range("B5").Value = "D:/Red-Nord/DISC T/Analiza"
